$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting (avoid numeric auto-conversion) for Price/Volume columns
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '42.317.07'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '2.231.88'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '244.50'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '0.622'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("D7").Value = '74.35'
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.619'
$ws.Range("E9").Value = '  +3.03%  '
$ws.Range("D10").Value = '43.48'
$ws.Range("E10").Value = '  +9.68%  '
$ws.Range("D11").Value = '0.0971'
$ws.Range("E11").Value = '  +4.35%  '
$ws.Range("D12").Value = '7.20'
$ws.Range("E12").Value = '  +3.10%  '
$ws.Range("E13").Value = '  +1.30%  '
$ws.Range("D14").Value = '14.39'
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").Value = '0.852'
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").Value = '2.266.53'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '42.173.26'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0000114'
$ws.Range("E18").Value = '  +18.79%  '
$ws.Range("D19").Value = '6.19'
$ws.Range("E19").Value = '  +3.43%  '
$ws.Range("D20").Value = '72.26'
$ws.Range("E20").Value = '  +1.42%  '
$ws.Range("D21").Value = '10.04'
$ws.Range("E21").Value = '  +41.04%  '
$ws.Range("D22").Value = '230.32'
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").Value = '2.19'
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("D24").Value = '11.88'
$ws.Range("E24").Value = '  +9.05%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").Value = '3.62'
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("E27").Value = '  +1.88%  '
$ws.Range("E28").Value = '  +3.65%  '
$ws.Range("D29").Value = '166.86'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").Value = '20.86'
$ws.Range("E30").Value = '  +2.53%  '
$ws.Range("D31").Value = '5.66'
$ws.Range("E31").Value = '  +18.70%  '
$ws.Range("D32").Value = '0.0807'
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").Value = '0.119'
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '29.88'
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.125'
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").Value = '4.42'
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("D37").Value = '0.0307'
$ws.Range("E37").Value = '  +4.49%  '
$ws.Range("D38").Value = '13.27'
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").Value = '2.17'
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("D40").Value = '5.64'
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("D41").Value = '64.06'
$ws.Range("E41").Value = '  +8.80%  '
$ws.Range("D42").Value = '0.202'
$ws.Range("E42").Value = '  +1.83%  '
$ws.Range("D43").Value = '8.84'
$ws.Range("E43").Value = '  +2.99%  '
$ws.Range("D44").Value = '105.57'
$ws.Range("E44").Value = '  -2.01%  '
$ws.Range("E45").Value = '  +3.02%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '2.43'
$ws.Range("E46").Value = '  +10.41%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '1.13'
$ws.Range("E47").Value = '  +3.10%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").Value = '1.18'
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").Value = '2.72'
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D50").Value = '4.09'
$ws.Range("E50").Value = '  +1.62%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.444.24'
$ws.Range("E51").Value = '  +0.75%  '
